$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price cells whose new values are plain decimal numbers,
# so Excel does not auto-convert them to numeric (losing trailing zeros / type).
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = "24.585.37"
$ws.Range("E2").Value = "  +3.86%  "
$ws.Range("D3").Value = "1.694.02"
$ws.Range("E3").Value = "  +2.30%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "316.02"
$ws.Range("E5").Value = "  +2.03%  "
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").Value = "0.3935"
$ws.Range("E7").Value = "  +1.39%  "
$ws.Range("D8").Value = "0.4013"
$ws.Range("E8").Value = "  +1.98%  "
$ws.Range("E9").Value = "  +7.38%  "
$ws.Range("B10").Value = "OKB"
$ws.Range("C10").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D10").Value = "54.45"
$ws.Range("E10").Value = "  +10.91%  "
$ws.Range("B11").Value = "BinanceUSD"
$ws.Range("C11").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D11").Value = "1.002"
$ws.Range("E11").Value = "  -0.11%  "
$ws.Range("D12").Value = "0.08758"
$ws.Range("E12").Value = "  +1.29%  "
$ws.Range("D13").Value = "7.197"
$ws.Range("E13").Value = "  +7.76%  "
$ws.Range("D14").Value = "23.20"
$ws.Range("E14").Value = "  +2.69%  "
$ws.Range("E15").Value = "  +0.72%  "
$ws.Range("D16").Value = "7.595"
$ws.Range("E16").Value = "  +4.99%  "
$ws.Range("D17").Value = "1.699.15"
$ws.Range("E17").Value = "  +2.29%  "
$ws.Range("D18").Value = "100.31"
$ws.Range("E18").Value = "  +0.71%  "
$ws.Range("D19").Value = "0.07048"
$ws.Range("E19").Value = "  +3.63%  "
$ws.Range("E20").Value = "  +3.28%  "
$ws.Range("D21").Value = "6.856"
$ws.Range("E21").Value = "  +2.86%  "
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("D23").Value = "14.01"
$ws.Range("E23").Value = "  +1.24%  "
$ws.Range("D24").Value = "24.576.37"
$ws.Range("E24").Value = "  +3.87%  "
$ws.Range("D25").Value = "3.008"
$ws.Range("E25").Value = "  +8.06%  "
$ws.Range("D26").Value = "2.310"
$ws.Range("E26").Value = "  -0.57%  "
$ws.Range("D27").Value = "22.34"
$ws.Range("E27").Value = "  +2.92%  "
$ws.Range("D28").Value = "159.11"
$ws.Range("E28").Value = "  +0.62%  "
$ws.Range("D29").Value = "5.207"
$ws.Range("E29").Value = "  +1.09%  "
$ws.Range("D30").Value = "133.90"
$ws.Range("E30").Value = "  +3.36%  "
$ws.Range("D31").Value = "7.457"
$ws.Range("E31").Value = "  +14.40%  "
$ws.Range("D32").Value = "1.883.80"
$ws.Range("E32").Value = "  +2.04%  "
$ws.Range("D33").Value = "1.093"
$ws.Range("E33").Value = "  -3.08%  "
$ws.Range("D34").Value = "7.281"
$ws.Range("E34").Value = "  +11.17%  "
$ws.Range("D35").Value = "0.08520"
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("D36").Value = "11.33"
$ws.Range("E36").Value = "  +9.72%  "
$ws.Range("D37").Value = "1.962"
$ws.Range("E37").Value = "  +0.96%  "
$ws.Range("D38").Value = "0.2721"
$ws.Range("D39").Value = "14.52"
$ws.Range("E39").Value = "  +0.79%  "
$ws.Range("D40").Value = "0.02744"
$ws.Range("E40").Value = "  +9.18%  "
$ws.Range("D41").Value = "0.09023"
$ws.Range("E41").Value = "  +2.85%  "
$ws.Range("D42").Value = "1.461"
$ws.Range("E42").Value = "  +0.92%  "
$ws.Range("D43").Value = "0.7660"
$ws.Range("E43").Value = "  +1.81%  "
$ws.Range("D44").Value = "0.7176"
$ws.Range("E44").Value = "  +2.62%  "
$ws.Range("D45").Value = "15.45"
$ws.Range("E45").Value = "  +4.18%  "
$ws.Range("D46").Value = "2.505"
$ws.Range("E46").Value = "  +4.45%  "
$ws.Range("D47").Value = "4.204"
$ws.Range("E47").Value = "  +2.83%  "
$ws.Range("E49").Value = "  +14.01%  "
$ws.Range("D50").Value = "141.17"
$ws.Range("E50").Value = "  +2.40%  "
$ws.Range("E51").Value = "  +3.05%  "
